# Case with 380 kV done: updated bus voltage magnitude results (vm_pu)
# after re-running the load-flow with the new 380 kV slack-bus setpoint.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letters B..N (skipping G/H, which are unaffected) mapped to indices
$colIndex = @{
    "B" = 2
    "C" = 3
    "D" = 4
    "E" = 5
    "F" = 6
    "I" = 9
    "J" = 10
    "K" = 11
    "L" = 12
    "M" = 13
    "N" = 14
}

# New vm_pu values per data row (row 2 = bus index 0 ... row 25 = bus index 23)
$newValues = @{
    2 = @{ "B" = 1.02; "C" = 1.044269252652559; "D" = 1.042574350693107; "E" = 1.05153667465833; "F" = 1.060888035034363; "I" = 1.041731533383036; "J" = 1.049335242836279; "K" = 1.045350653450992; "L" = 1.054287903491258; "M" = 1.063613590662148; "N" = 1.020252522136109 }
    3 = @{ "B" = 1.02; "C" = 1.045378982712692; "D" = 1.043170791527352; "E" = 1.052530256005684; "F" = 1.06198092762699; "I" = 1.041978314892545; "J" = 1.050091365997522; "K" = 1.045758515411779; "L" = 1.055093685577495; "M" = 1.064520294503491; "N" = 1.020509636479818 }
    4 = @{ "B" = 1.02; "C" = 1.046097198055176; "D" = 1.043556793396399; "E" = 1.053173647814641; "F" = 1.062688719804192; "I" = 1.042136875497875; "J" = 1.050580213184197; "K" = 1.046021796759023; "L" = 1.055614944087508; "M" = 1.065107013679338; "N" = 1.020675720670588 }
    5 = @{ "B" = 1.02; "C" = 1.046399171189773; "D" = 1.04371908287049; "E" = 1.053444244065106; "F" = 1.062986423412037; "I" = 1.042203265134253; "J" = 1.050785625169838; "K" = 1.046132328057453; "L" = 1.05583404844248; "M" = 1.065353675408081; "N" = 1.02074547388707 }
    6 = @{ "B" = 1.02; "C" = 1.046449875898689; "D" = 1.04374633278971; "E" = 1.053489685066318; "F" = 1.063036417821309; "I" = 1.042214396450134; "J" = 1.050820108918771; "K" = 1.046150877803636; "L" = 1.055870835106846; "M" = 1.065395091287923; "N" = 1.02075718173437 }
    7 = @{ "B" = 1.02; "C" = 1.046101232895222; "D" = 1.043558961861689; "E" = 1.053177263085922; "F" = 1.062692697153798; "I" = 1.042137763657713; "J" = 1.05058295829967; "K" = 1.046023274282527; "L" = 1.055617871899289; "M" = 1.065110309565929; "N" = 1.020676652986083 }
    8 = @{ "B" = 1.02; "C" = 1.044644261237206; "D" = 1.042775906525047; "E" = 1.051872361009552; "F" = 1.061257255280927; "I" = 1.041815166824514; "J" = 1.049590864135959; "K" = 1.045488622858504; "L" = 1.054560249514634; "M" = 1.063920011041392; "N" = 1.020339474279386 }
    9 = @{ "B" = 1.02; "C" = 1.042077984816745; "D" = 1.041396612497119; "E" = 1.049576626076656; "F" = 1.05873255619793; "I" = 1.041238116686148; "J" = 1.047839488456404; "K" = 1.044541685748655; "L" = 1.052695542050711; "M" = 1.061822718093432; "N" = 1.019743136850173 }
    10 = @{ "B" = 1.02; "C" = 1.040367834769296; "D" = 1.040477523940739; "E" = 1.048048609193641; "F" = 1.05705260652008; "I" = 1.040847653728597; "J" = 1.046669764604452; "K" = 1.043907200759933; "L" = 1.05145170471533; "M" = 1.060424635282511; "N" = 1.019344114842671 }
    11 = @{ "B" = 1.02; "C" = 1.039627478634734; "D" = 1.040079667795023; "E" = 1.04738754799368; "F" = 1.056325923954072; "I" = 1.040677214897556; "J" = 1.046162752284647; "K" = 1.043631712941083; "L" = 1.05091294309611; "M" = 1.05981927610321; "N" = 1.019170987654326 }
    12 = @{ "B" = 1.02; "C" = 1.039352499116617; "D" = 1.03993190479089; "E" = 1.047142087527325; "F" = 1.0560561137046; "I" = 1.040613701308892; "J" = 1.045974347869798; "K" = 1.043529272232974; "L" = 1.050712797082559; "M" = 1.059594421587615; "N" = 1.019106628259499 }
    13 = @{ "B" = 1.02; "C" = 1.039411482187505; "D" = 1.039963599601818; "E" = 1.047194735690394; "F" = 1.05611398385309; "I" = 1.040627334463909; "J" = 1.046014764765224; "K" = 1.043551251182342; "L" = 1.050755730258884; "M" = 1.059642653518552; "N" = 1.019120435931909 }
    14 = @{ "B" = 1.02; "C" = 1.039604748303854; "D" = 1.040067453275238; "E" = 1.047367256378609; "F" = 1.056303619070705; "I" = 1.040671969025782; "J" = 1.046147180298438; "K" = 1.04362324744659; "L" = 1.050896399494582; "M" = 1.059800689501808; "N" = 1.019165668752307 }
    15 = @{ "B" = 1.02; "C" = 1.039723828789082; "D" = 1.040131443434982; "E" = 1.047473563607635; "F" = 1.056420474424333; "I" = 1.040699442690912; "J" = 1.046228755607481; "K" = 1.04366759186447; "L" = 1.050983067009214; "M" = 1.059898061058347; "N" = 1.019193531263524 }
    16 = @{ "B" = 1.02; "C" = 1.040416972858793; "D" = 1.040503930870941; "E" = 1.048092493866312; "F" = 1.057100849766347; "I" = 1.040858936431874; "J" = 1.046703402490822; "K" = 1.043925468199379; "L" = 1.051487456916904; "M" = 1.060464811383786; "N" = 1.019355597394757 }
    17 = @{ "B" = 1.02; "C" = 1.040851803419479; "D" = 1.040737614082757; "E" = 1.048480887750236; "F" = 1.057527831039058; "I" = 1.040958617104824; "J" = 1.047000998129444; "K" = 1.044087026443442; "L" = 1.051803801272015; "M" = 1.060820324110084; "N" = 1.01945716404798 }
    18 = @{ "B" = 1.02; "C" = 1.041105447269433; "D" = 1.040873928590052; "E" = 1.048707487316621; "F" = 1.057776954051932; "I" = 1.041016627363683; "J" = 1.047174531034721; "K" = 1.044181188068284; "L" = 1.051988303011167; "M" = 1.061027690599162; "N" = 1.019516372576428 }
    19 = @{ "B" = 1.02; "C" = 1.041191935759451; "D" = 1.040920410151969; "E" = 1.048784761411075; "F" = 1.057861910789503; "I" = 1.041036385003384; "J" = 1.047233692865202; "K" = 1.044213282420565; "L" = 1.05205121052986; "M" = 1.061098397545251; "N" = 1.019536555457278 }
    20 = @{ "B" = 1.02; "C" = 1.040805148717402; "D" = 1.040712540934953; "E" = 1.048439210971815; "F" = 1.057482012537003; "I" = 1.040947935942452; "J" = 1.046969074065214; "K" = 1.044069700273833; "L" = 1.051769862231505; "M" = 1.060782180760386; "N" = 1.019446270380538 }
    21 = @{ "B" = 1.02; "C" = 1.039547835656915; "D" = 1.040036870437345; "E" = 1.047316450932152; "F" = 1.056247773135008; "I" = 1.040658830917872; "J" = 1.046108189364405; "K" = 1.043602049435049; "L" = 1.050854976628511; "M" = 1.059754151766583; "N" = 1.019152350256108 }
    22 = @{ "B" = 1.02; "C" = 1.038757438254369; "D" = 1.039612157087938; "E" = 1.046611031015438; "F" = 1.055472405811495; "I" = 1.040475873031591; "J" = 1.045566468692769; "K" = 1.043307369578843; "L" = 1.050279601473276; "M" = 1.059107805476416; "N" = 1.018967248732804 }
    23 = @{ "B" = 1.02; "C" = 1.039176431282201; "D" = 1.039837295091779; "E" = 1.046984939721442; "F" = 1.055883381338205; "I" = 1.040572974836403; "J" = 1.045853687635172; "K" = 1.043463646264891; "L" = 1.050584633020105; "M" = 1.059450444376009; "N" = 1.019065403192046 }
    24 = @{ "B" = 1.02; "C" = 1.040826229923815; "D" = 1.040723870377172; "E" = 1.0484580427409; "F" = 1.05750271572302; "I" = 1.040952762707149; "J" = 1.046983499328769; "K" = 1.044077529448337; "L" = 1.051785197874079; "M" = 1.06079941609213; "N" = 1.019451192863669 }
    25 = @{ "B" = 1.02; "C" = 1.042741302729563; "D" = 1.041753120574562; "E" = 1.050169693134977; "F" = 1.059384690464629; "I" = 1.041388314807907; "J" = 1.048292638057202; "K" = 1.044787057675975; "L" = 1.053177736935638; "M" = 1.062364898875672; "N" = 1.019897562773423 }
}

foreach ($row in $newValues.Keys) {
    $rowData = $newValues[$row]
    foreach ($col in $rowData.Keys) {
        $ws.Cells.Item([int]$row, $colIndex[$col]).Value = $rowData[$col]
    }
}
